$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Turn the URL inside the <dryingUrl> line into a real hyperlink.
#    Original text:
#      <dryingUrl>https://apps.chhs.colostate.edu/preservesmart/produce/drying/brussels-sprouts/</dryingUrl>
#    The hyperlinked span covers everything from just after "<dryingUrl>"
#    up to (but not including) the final ">" character - i.e. it swallows
#    the literal "</dryingUrl" text too, leaving a lone "<dryingUrl>" run
#    before it and a lone ">" run after it.
# ------------------------------------------------------------------
$dryingPara = $d.Paragraphs(3)
$paraRange = $dryingPara.Range

$openTag = "<dryingUrl>"
$hyperStart = $paraRange.Start + $openTag.Length
$hyperEnd = $paraRange.End - 2   # -1 for the paragraph mark, -1 more to exclude the trailing ">"

$hyperRange = $d.Range($hyperStart, $hyperEnd)

# Apply the "Hyperlink" character style to the run first so the style
# definition gets written into styles.xml ...
$hyperRange.Style = "Hyperlink"

# ... then turn that same span into an actual hyperlink field.
$hyperRange = $d.Range($hyperStart, $hyperEnd)
$linkAddress = "https://apps.chhs.colostate.edu/preservesmart/produce/drying/brussels-sprouts/"
$d.Hyperlinks.Add($hyperRange, $linkAddress) | Out-Null

$hyperlinkStyle = $d.Styles("Hyperlink")
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.UnhideWhenUsed = $true

# ------------------------------------------------------------------
# 2) After the <dryingUrl> paragraph, add two blank paragraphs and a
#    new paragraph containing the picture link.
# ------------------------------------------------------------------
$dryingPara = $d.Paragraphs(3)
$dryingPara.Range.InsertParagraphAfter() | Out-Null

$blank1 = $d.Paragraphs(4)
$blank1.Range.InsertParagraphAfter() | Out-Null

$blank2 = $d.Paragraphs(5)
$blank2.Range.InsertParagraphAfter() | Out-Null

$picturePara = $d.Paragraphs(6)
$picturePara.Range.Text = "https://i.imgur.com/gRLAR86.jpg"

# ------------------------------------------------------------------
# 3) Register the (otherwise unused) "Unresolved Mention" character
#    style that Word adds alongside "Hyperlink" in modern documents.
# ------------------------------------------------------------------
$mentionStyle = $d.Styles.Add("UnresolvedMention", 2)
$mentionStyle.NameLocal = "Unresolved Mention"
$mentionStyle.BaseStyle = "DefaultParagraphFont"
$mentionStyle.Priority = 99
$mentionStyle.UnhideWhenUsed = $true
$mentionStyle.Font.Color = 6053472
